$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$aValues = @("g4","g4","e4","a4","g4","c5","g4","d4","e4","d4","d4","g3","c4","d4","e4","d4","e4","f4","e4","g3","c4","d4","e4","f4","e4","c4","c4","g4","f4","e4","g3","g3","g4","g4","a4","d4","d4","g3","c4","d4","e4","g3","c4","d4","e4","d4","e4","f4","g4","c4","c4","g4","g4","g4","a4","c5","c5","d5","e5","d5","c5","a4","d5","c5","a4","d5","c5","g4","e4","a4","g4","c5","g4","e4","a4","g4","c5","g4","d4","e4","f4","g4","c4","c4","g3","b3","c4","c5","g4")
$bValues = @(2,4,4,4,4,4,8,8,4,8,2,8,8,8,8,8,8,8,8,8,4,8,2,8,8,8,8,8,8,8,4,8,2,8,8,8,8,8,8,8,8,8,8,8,8,4,8,4,8,8,8,8,8,8,8,4,8,2,8,8,8,8,8,8,8,8,8,8,8,8,8,4,4,4,8,2,8,8,8,8,8,8,4,8,4,8,4,2,4)

$n = 89

# Build a 2D array (n x 2) for bulk assignment to A1:B89
$data = New-Object 'object[,]' $n,2
for ($i = 0; $i -lt $n; $i++) {
    $data[$i,0] = $aValues[$i]
    $data[$i,1] = $bValues[$i]
}

$range = $ws.Range("A1:B$n")
$range.Value = $data
